$wb = $excel.ActiveWorkbook

# --- Remove the extra empty "Sheet" tab created by mistake when the workbook
#     was first generated (bugfix: no extra sheet on new workbook creation). ---
$excel.DisplayAlerts = $false
foreach ($s in @($wb.Worksheets)) {
    if ($s.Name -eq "Sheet") {
        $s.Delete()
    }
}

$ws = $wb.Worksheets("Car Evaluations")

# --- Row 2: only the Update Date timestamp changed (re-scraped at a later time) ---
$ws.Range("O2").Value = "2019-12-29 19:32:08.884654"

# --- Row 3: replaced with freshly scraped listing data ---
$ws.Range("A3").Value = "2017 Honda Accord LX"
$ws.Range("B3").Value = "'$15,995"
$ws.Range("B3").ClearFormats()
$ws.Range("C3").Value = "26,832 miles "
$ws.Range("D3").Value = "Falls Church, VA "
$ws.Range("E3").Value = "Unspecified"
$ws.Range("F3").Value = "Unspecified"
$ws.Range("G3").Value = "FWD"
$ws.Range("H3").Value = "Automatic"
$ws.Range("I3").Value = "Sedan"
$ws.Range("J3").Value = "4 Cyl 2.4 L"
$ws.Range("K3").Value = "Gasoline"
$ws.Range("L3").Value = "26/34"
$ws.Range("M3").Value = "1HGCR2F35HA094805"
$ws.Range("N3").Value = "HP4847"
$ws.Range("O3").Value = "2019-12-29 19:33:12.889521"
$ws.Range("P3").Value = "https://www.carfax.com/vehicle/1HGCR2F35HA094805"

# --- Row 4: replaced with freshly scraped listing data ---
$ws.Range("A4").Value = "2016 Honda Accord Sport"
$ws.Range("B4").Value = "'$16,277"
$ws.Range("B4").ClearFormats()
$ws.Range("C4").Value = "28,474 miles "
$ws.Range("D4").Value = "Stafford, VA "
$ws.Range("E4").Value = "Red"
$ws.Range("F4").Value = "Black"
$ws.Range("G4").Value = "FWD"
$ws.Range("H4").Value = "Automatic"
$ws.Range("I4").Value = "Sedan"
$ws.Range("J4").Value = "4 Cyl 2.4 L"
$ws.Range("K4").Value = "Gasoline"
$ws.Range("L4").Value = "26/34"
$ws.Range("M4").Value = "1HGCR2F52GA125671"
$ws.Range("N4").Value = "'125671"
$ws.Range("N4").ClearFormats()
$ws.Range("O4").Value = "2019-12-29 19:33:14.067253"
$ws.Range("P4").Value = "https://www.carfax.com/vehicle/1HGCR2F52GA125671"

# --- Row 5: brand new listing appended to the bottom of the table ---
$ws.Range("A5").Value = "2015 Honda Accord EXL"
$ws.Range("B5").Value = "'$15,371"
$ws.Range("B5").ClearFormats()
$ws.Range("C5").Value = "68,777 miles "
$ws.Range("D5").Value = "Germantown, MD "
$ws.Range("E5").Value = "Green"
$ws.Range("F5").Value = "Black"
$ws.Range("G5").Value = "FWD"
$ws.Range("H5").Value = "Automatic"
$ws.Range("I5").Value = "Sedan"
$ws.Range("J5").Value = "6 Cyl 3.5 L"
$ws.Range("K5").Value = "Gasoline"
$ws.Range("L5").Value = "21/31"
$ws.Range("M5").Value = "1HGCR3F85FA011615"
$ws.Range("N5").Value = "H191579B"
$ws.Range("O5").Value = "2019-12-29 19:34:34.986326"
$ws.Range("P5").Value = "https://www.carfax.com/vehicle/1HGCR3F85FA011615"
